$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E and F, rows 2-7 hold fractional rates that need to become
# percentage points (i.e. multiplied by 100), per the ranking table update.
for ($row = 2; $row -le 7; $row++) {
    foreach ($col in 'E', 'F') {
        $cell = $ws.Range("$col$row")
        $current = $cell.Value()
        $cell.Value = $current * 100
    }
}
